$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing B-column values (resampled error series) ---
$ws.Range("B8").Value  = -1.164515107544459
$ws.Range("B9").Value  = -0.9686002503882047
$ws.Range("B13").Value = 0.6076553625472806
$ws.Range("B16").Value = 0.4380760559974082
$ws.Range("B18").Value = -0.4191283315245076
$ws.Range("B19").Value = 1.035245970010962
$ws.Range("B20").Value = 0.2459967371659499
$ws.Range("B21").Value = 0.9978217259122815
$ws.Range("B22").Value = -0.3737868047750048
$ws.Range("B23").Value = 0.2315426864241067

# --- Append two new rows (2025-07-01_diff, 2025-10-01_diff) ---
$ws.Range("A23").Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)

$ws.Range("A24").Value = "2025-07-01_diff"
$ws.Range("A25").Value = "2025-10-01_diff"

$excel.CutCopyMode = 0
